$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25: add new values for columns J..O ---
$ws.Range("J25").Value = 0.041666666666666664
$ws.Range("J25").NumberFormat = "h:mm"

$ws.Range("K25").Value = 0.1111111111111111
$ws.Range("K25").NumberFormat = "h:mm"

$ws.Range("L25").Value = 0.09375
$ws.Range("L25").NumberFormat = "h:mm"

$ws.Range("M25").Value = 0.15277777777777776
$ws.Range("M25").NumberFormat = "h:mm"

$ws.Range("N25").Value = 0.15625
$ws.Range("N25").NumberFormat = "h:mm"

$ws.Range("O25").Value = 0.06944444444444443
$ws.Range("O25").NumberFormat = "h:mm"

# --- New row 31 ---
$ws.Range("A31").Value = "Faire la page rechercher"
$ws.Range("J31").Value = 0.020833333333333332
$ws.Range("J31").NumberFormat = "h:mm"

# --- New row 32 ---
$ws.Range("A32").Value = "Créer la vue dans la BD"
$ws.Range("J32").Value = 0.020833333333333332
$ws.Range("J32").NumberFormat = "h:mm"
$ws.Range("K32").NumberFormat = "h:mm"

# --- New row 33 ---
$ws.Range("A33").Value = "Fonction et page rechercher"
$ws.Range("J33").Value = 0.07291666666666667
$ws.Range("J33").NumberFormat = "h:mm"
$ws.Range("K33").Value = 0.041666666666666664
$ws.Range("K33").NumberFormat = "h:mm"
$ws.Range("L33").Value = 0.0625
$ws.Range("L33").NumberFormat = "h:mm"
$ws.Range("O33").Value = 0.08333333333333333
$ws.Range("O33").NumberFormat = "h:mm"

# --- Update selection to match target ---
$ws.Range("O26").Select()
